$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-11 20:18:40"
$ws.Range("H2").Formula = '=TEXT("81%","@")'
$ws.Range("H2").Copy()
$ws.Range("H2").PasteSpecial(-4163)
$ws.Range("I2").Value = "8.7 mm"
$ws.Range("N2").Value = "0.4 °C 19:59 TU"
$ws.Range("O2").Value = "3.0 °C"
$ws.Range("E3").Value = "2026-02-11 20:18:42"
$ws.Range("I3").Value = "3.9 mm"
$ws.Range("E4").Value = "2026-02-11 20:18:45"
$ws.Range("I4").Value = "0.2 mm"
$ws.Range("J4").Value = "1002.2 hPa"
$ws.Range("E5").Value = "2026-02-11 20:18:47"
$ws.Range("I5").Value = "3.7 mm"
$ws.Range("E6").Value = "2026-02-11 20:18:50"
$ws.Range("J6").Value = "1002.7 hPa"
$ws.Range("E7").Value = "2026-02-11 20:18:53"
$ws.Range("H7").Formula = '=TEXT("44%","@")'
$ws.Range("H7").Copy()
$ws.Range("H7").PasteSpecial(-4163)
$ws.Range("I7").Value = "0.4 mm"
$ws.Range("N7").Value = "15.2 °C 19:59 TU"
$ws.Range("O7").Value = "19.0 °C"
$ws.Range("E8").Value = "2026-02-11 20:18:55"
$ws.Range("H8").Formula = '=TEXT("56%","@")'
$ws.Range("H8").Copy()
$ws.Range("H8").PasteSpecial(-4163)
$ws.Range("I8").Value = "1.5 mm"
$ws.Range("N8").Value = "11.2 °C 19:59 TU"
$ws.Range("O8").Value = "15.1 °C"
$ws.Range("E9").Value = "2026-02-11 20:18:57"
$ws.Range("E10").Value = "2026-02-11 20:19:00"
$ws.Range("H10").Formula = '=TEXT("74%","@")'
$ws.Range("H10").Copy()
$ws.Range("H10").PasteSpecial(-4163)
$ws.Range("I10").Value = "0.1 mm"
$ws.Range("O10").Value = "13.8 °C"
$ws.Range("E11").Value = "2026-02-11 20:19:02"
$ws.Range("I11").Value = "1.2 mm"
$ws.Range("E12").Value = "2026-02-11 20:19:04"
$ws.Range("E13").Value = "2026-02-11 20:19:06"
$ws.Range("I13").Value = "1.1 mm"
$ws.Range("J13").Value = "1005.1 hPa"
$ws.Range("E14").Value = "2026-02-11 20:19:09"
$ws.Range("H14").Formula = '=TEXT("48%","@")'
$ws.Range("H14").Copy()
$ws.Range("H14").PasteSpecial(-4163)
$ws.Range("N14").Value = "14.1 °C 19:58 TU"
$ws.Range("O14").Value = "18.9 °C"
$ws.Range("E15").Value = "2026-02-11 20:19:12"
$ws.Range("E16").Value = "2026-02-11 20:19:14"
$ws.Range("G16").Value = "79 cm"
$ws.Range("I16").Value = "8.2 mm"
$ws.Range("E17").Value = "2026-02-11 20:19:17"
$ws.Range("I17").Value = "5.3 mm"
$ws.Range("E18").Value = "2026-02-11 20:19:19"
$ws.Range("H18").Formula = '=TEXT("69%","@")'
$ws.Range("H18").Copy()
$ws.Range("H18").PasteSpecial(-4163)
$ws.Range("J18").Value = "1002.8 hPa"
$ws.Range("L18").Value = "40.7 km/h - 228º 19:40 TU"
$ws.Range("O18").Value = "14.3 °C"
$ws.Range("E19").Value = "2026-02-11 20:19:22"
$ws.Range("I19").Value = "1.5 mm"
$ws.Range("E20").Value = "2026-02-11 20:19:25"
$ws.Range("I20").Value = "1.8 mm"
$ws.Range("O20").Value = "-1.0 °C"
$ws.Range("E21").Value = "2026-02-11 20:19:27"
$ws.Range("I21").Value = "3.9 mm"
$ws.Range("J21").Value = "1005.6 hPa"
$ws.Range("O21").Value = "8.3 °C"
$ws.Range("E22").Value = "2026-02-11 20:19:30"
$ws.Range("I22").Value = "1.9 mm"
$ws.Range("M22").Value = "-0.1 °C 19:42 TU"
$ws.Range("E23").Value = "2026-02-11 20:19:32"
$ws.Range("H23").Formula = '=TEXT("74%","@")'
$ws.Range("H23").Copy()
$ws.Range("H23").PasteSpecial(-4163)
$ws.Range("I23").Value = "6.9 mm"
$ws.Range("E24").Value = "2026-02-11 20:19:35"
$ws.Range("I24").Value = "8.6 mm"
$ws.Range("J24").Value = "1006.8 hPa"
$ws.Range("N24").Value = "10.9 °C 19:59 TU"
$ws.Range("O24").Value = "13.2 °C"
$ws.Range("E25").Value = "2026-02-11 20:19:38"
$ws.Range("I25").Value = "2.9 mm"
$ws.Range("O25").Value = "1.7 °C"
$ws.Range("E26").Value = "2026-02-11 20:19:40"
$ws.Range("J26").Value = "1002.6 hPa"
$ws.Range("E27").Value = "2026-02-11 20:19:43"
$ws.Range("I27").Value = "2.4 mm"
$ws.Range("E28").Value = "2026-02-11 20:19:46"
$ws.Range("H28").Formula = '=TEXT("79%","@")'
$ws.Range("H28").Copy()
$ws.Range("H28").PasteSpecial(-4163)
$ws.Range("J28").Value = "1002.9 hPa"
$ws.Range("L28").Value = "51.5 km/h - 281º 19:58 TU"
$ws.Range("M28").Value = "16.8 °C 19:46 TU"
$ws.Range("O28").Value = "11.2 °C"
$ws.Range("E29").Value = "2026-02-11 20:19:49"
$ws.Range("H29").Formula = '=TEXT("84%","@")'
$ws.Range("H29").Copy()
$ws.Range("H29").PasteSpecial(-4163)
$ws.Range("L29").Value = "56.2 km/h - 244º 19:39 TU"
$ws.Range("O29").Value = "13.0 °C"
$ws.Range("E30").Value = "2026-02-11 20:19:52"
$ws.Range("J30").Value = "1002.9 hPa"
$ws.Range("E31").Value = "2026-02-11 20:19:54"
$ws.Range("J31").Value = "1002.1 hPa"
$ws.Range("K31").Value = "9.4 MJ/m2"
$ws.Range("E32").Value = "2026-02-11 20:19:57"
$ws.Range("I32").Value = "4.1 mm"
$ws.Range("E33").Value = "2026-02-11 20:20:00"
$ws.Range("H33").Formula = '=TEXT("81%","@")'
$ws.Range("H33").Copy()
$ws.Range("H33").PasteSpecial(-4163)
$ws.Range("I33").Value = "2.3 mm"
$ws.Range("J33").Value = "1004.7 hPa"
$ws.Range("E34").Value = "2026-02-11 20:20:02"
$ws.Range("I34").Value = "3.0 mm"
$ws.Range("E35").Value = "2026-02-11 20:20:05"
$ws.Range("J35").Value = "1007.3 hPa"
$ws.Range("K35").Value = "7.2 MJ/m2"
$ws.Range("M35").Value = "14.5 °C 19:39 TU"
$ws.Range("O35").Value = "11.1 °C"
$ws.Range("E36").Value = "2026-02-11 20:20:07"
$ws.Range("J36").Value = "1003.0 hPa"
$ws.Range("E37").Value = "2026-02-11 20:20:10"
$ws.Range("I37").Value = "0.9 mm"
$ws.Range("J37").Value = "1004.2 hPa"
$ws.Range("E38").Value = "2026-02-11 20:20:13"
$ws.Range("H38").Formula = '=TEXT("59%","@")'
$ws.Range("H38").Copy()
$ws.Range("H38").PasteSpecial(-4163)
$ws.Range("I38").Value = "2.3 mm"
$ws.Range("L38").Value = "74.2 km/h - 255º 19:52 TU"
$ws.Range("E39").Value = "2026-02-11 20:20:15"
$ws.Range("H39").Formula = '=TEXT("57%","@")'
$ws.Range("H39").Copy()
$ws.Range("H39").PasteSpecial(-4163)
$ws.Range("E40").Value = "2026-02-11 20:20:18"
$ws.Range("H40").Formula = '=TEXT("91%","@")'
$ws.Range("H40").Copy()
$ws.Range("H40").PasteSpecial(-4163)
$ws.Range("I40").Value = "5.6 mm"
$ws.Range("J40").Value = "1006.8 hPa"
$ws.Range("E41").Value = "2026-02-11 20:20:21"
$ws.Range("H41").Formula = '=TEXT("49%","@")'
$ws.Range("H41").Copy()
$ws.Range("H41").PasteSpecial(-4163)
$ws.Range("J41").Value = "1004.6 hPa"
$ws.Range("N41").Value = "13.6 °C 19:31 TU"
$ws.Range("O41").Value = "18.8 °C"
$ws.Range("E42").Value = "2026-02-11 20:20:24"
$ws.Range("H42").Formula = '=TEXT("89%","@")'
$ws.Range("H42").Copy()
$ws.Range("H42").PasteSpecial(-4163)
$ws.Range("O42").Value = "12.7 °C"
$ws.Range("E43").Value = "2026-02-11 20:20:26"
$ws.Range("H43").Formula = '=TEXT("64%","@")'
$ws.Range("H43").Copy()
$ws.Range("H43").PasteSpecial(-4163)
$ws.Range("I43").Value = "6.1 mm"
$ws.Range("N43").Value = "9.8 °C 19:40 TU"
$ws.Range("O43").Value = "13.1 °C"
$ws.Range("E44").Value = "2026-02-11 20:20:29"
$ws.Range("G44").Value = "220 cm"
$ws.Range("I44").Value = "6.9 mm"
$ws.Range("E45").Value = "2026-02-11 20:20:32"
$ws.Range("I45").Value = "4.9 mm"
$ws.Range("J45").Value = "1005.6 hPa"
$ws.Range("E46").Value = "2026-02-11 20:20:34"
$ws.Range("H46").Formula = '=TEXT("61%","@")'
$ws.Range("H46").Copy()
$ws.Range("H46").PasteSpecial(-4163)
$ws.Range("J46").Value = "1007.1 hPa"
$ws.Range("N46").Value = "12.7 °C 19:46 TU"
$ws.Range("O46").Value = "16.9 °C"

$excel.CutCopyMode = $false

